$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows A8:G12 hold a small frequency table whose data (including the row
# labels in column A) need to be cyclically shifted up by one row: row 9's
# content moves into row 8, row 10's into row 9, row 11's into row 10,
# row 12's into row 11, and the old row 8's content wraps around into the
# new row 12. Capture the five rows first so the shift doesn't clobber
# values we still need to read.

$rows = @(8, 9, 10, 11, 12)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2
    )
}

# new row r (8..11) <= old row r+1 ; new row 12 <= old row 8
$order = @(9, 10, 11, 12, 8)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $destRow = $rows[$i]
    $srcRow = $order[$i]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le 7; $c++) {
        $v = $vals[$c - 1]
        $isEmpty = ($null -eq $v) -or (($v -is [string]) -and ($v.Length -eq 0))
        if ($isEmpty) {
            $ws.Cells.Item($destRow, $c).ClearContents()
        } else {
            $ws.Cells.Item($destRow, $c).Value = $v
        }
    }
}

# The label that used to sit in row 8 (an empty string) now lands in row 12's
# column A; it should end up as a genuinely empty cell, not a cell holding "".
$ws.Cells.Item(12, 1).ClearContents()
